# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New strikeout (K) values for column G, rows 2-37 (replacing old Strike# derived values)
$kValues = @{
    2  = 6
    3  = 1
    4  = 2
    5  = 3
    6  = 0
    7  = 4
    8  = 0
    9  = 0
    10 = 1
    11 = 5
    12 = 9
    13 = 7
    14 = 7
    15 = 2
    16 = 6
    17 = 3
    18 = 6
    19 = 4
    20 = 7
    21 = 3
    22 = 5
    23 = 6
    24 = 0
    25 = 6
    26 = 4
    27 = 7
    28 = 4
    29 = 4
    30 = 3
    31 = 6
    32 = 3
    33 = 8
    34 = 4
    35 = 6
    36 = 3
    37 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
